$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.019.36'
$ws.Range("E2").Value = '  +0.34%  '

$ws.Range("D3").Value = '3.086.05'
$ws.Range("E3").Value = '  +0.39%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.71'
$ws.Range("E5").Value = '  -0.93%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.00'
$ws.Range("E6").Value = '  +5.45%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").Value = '3.080.58'
$ws.Range("E8").Value = '  +0.24%  '

$ws.Range("E9").Value = '  -0.25%  '

$ws.Range("E10").Value = '  -0.02%  '

$ws.Range("E11").Value = '  +1.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.466'
$ws.Range("E12").Value = '  -1.07%  '

$ws.Range("E13").Value = '  +0.14%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.83'
$ws.Range("E14").Value = '  -0.76%  '

$ws.Range("E15").Value = '  +0.94%  '

$ws.Range("D16").Value = '3.603.65'
$ws.Range("E16").Value = '  +0.56%  '

$ws.Range("D17").Value = '67.015.29'
$ws.Range("E17").Value = '  +0.45%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.00'
$ws.Range("E18").Value = '  +0.10%  '

$ws.Range("D19").Value = '3.085.07'
$ws.Range("E19").Value = '  +0.34%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.44'
$ws.Range("E20").Value = '  +0.24%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '486.19'
$ws.Range("E21").Value = '  +0.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.67'
$ws.Range("E22").Value = '  -0.70%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.683'
$ws.Range("E23").Value = '  -0.30%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.28'
$ws.Range("E24").Value = '  +0.91%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.65'
$ws.Range("E25").Value = '  -1.11%  '

$ws.Range("E26").Value = '  +1.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.17'
$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("E28").Value = '  +0.06%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.84'
$ws.Range("E29").Value = '  +1.68%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.28'
$ws.Range("E30").Value = '  -0.99%  '

$ws.Range("E31").Value = '  -1.62%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.95'
$ws.Range("E32").Value = '  +0.77%  '

$ws.Range("E33").Value = '  +0.44%  '

$ws.Range("D34").Value = '0.0₃0942'
$ws.Range("E34").Value = '  +4.27%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.09%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '47.31'
$ws.Range("E36").Value = '  +2.54%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.311'
$ws.Range("E39").Value = '  +3.37%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '48.98'
$ws.Range("E40").Value = '  -0.60%  '

$ws.Range("E41").Value = '  +1.48%  '

$ws.Range("E42").Value = '  +0.32%  '

$ws.Range("E43").Value = '  +11.40%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.21'
$ws.Range("E44").Value = '  -1.09%  '

$ws.Range("D45").Value = '2.790.78'
$ws.Range("E45").Value = '  +0.93%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '367.91'
$ws.Range("E46").Value = '  -0.57%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0343'
$ws.Range("E47").Value = '  -0.32%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '134.27'
$ws.Range("E48").Value = '  -0.82%  '

$ws.Range("E49").Value = '  +0.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.49'
$ws.Range("E50").Value = '  +4.89%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.31'
$ws.Range("E51").Value = '  +7.95%  '

# Row 37 (was Mantle) now Filecoin
$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.55'
$ws.Range("E37").Value = '  -2.31%  '

# Row 38 (was Filecoin) now Mantle
$ws.Range("B38").Value = 'Mantle'
$ws.Range("C38").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.939'
$ws.Range("E38").Value = '  -1.37%  '
